# "Colocando header nos gráficos" - add a header label to column A on each
# data sheet (used by the charts), fix missing Portuguese accents, drop the
# now-unused "Teto" row from the emissions sheet, and refresh the cost
# sheet's header/values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: stamp a header into A1 using the same look as the existing
# header row (B1:E1 use style index 1 - bold / bordered / centered) by
# copying B1's format onto A1, then writing the label text.
# ---------------------------------------------------------------------
function Set-HeaderCell {
    param($ws, [string]$text)
    $ws.Range("B1").Copy() | Out-Null
    $ws.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("A1").Value = $text
}

# Sheets 1-4 share the same "Fonte/Tecnologia" row layout (rows 2-12).
$fonteSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($sheetName in $fonteSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    Set-HeaderCell $ws "Fonte/Tecnologia"

    # Row labels lose their bold/border style (s="1" -> no style) and a
    # few get their accents restored / abbreviations tidied up.
    $ws.Range("A2").ClearFormats()
    $ws.Range("A2").Value = "Hidro"

    $ws.Range("A3").ClearFormats()
    $ws.Range("A3").Value = "Gás Natural"

    $ws.Range("A4").ClearFormats()
    $ws.Range("A4").Value = "Carvão"

    $ws.Range("A5").ClearFormats()
    $ws.Range("A5").Value = "Nuclear"

    $ws.Range("A6").ClearFormats()
    $ws.Range("A6").Value = "Óleos Comb"

    $ws.Range("A7").ClearFormats()
    $ws.Range("A7").Value = "Biomassa"

    $ws.Range("A8").ClearFormats()
    $ws.Range("A8").Value = "Eólica"

    $ws.Range("A9").ClearFormats()
    $ws.Range("A9").Value = "Solar"

    $ws.Range("A10").ClearFormats()
    $ws.Range("A10").Value = "Outros"

    $ws.Range("A11").ClearFormats()
    $ws.Range("A11").Value = "Pot. Compl."

    $ws.Range("A12").ClearFormats()
    $ws.Range("A12").Value = "GD"
}

# ---------------------------------------------------------------------
# Sheet 5: Emissoes Totais (MtCO2eq)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")

Set-HeaderCell $ws "Período"

$ws.Range("A2").ClearFormats()
$ws.Range("A2").Value = "P.Médio"

$ws.Range("A3").ClearFormats()
$ws.Range("A3").Value = "P.Crítico"

# Drop the obsolete "Teto" row entirely (shifts dimension to A1:E3).
$ws.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6: Custo Total (bilhões de R$)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Custo Total (bilhões de R$)")

Set-HeaderCell $ws "Tipo Expansão"

# B1 used to read "Custo"; it now carries the year label, kept as text
# (not a number) so it still reads "2015" rather than becoming numeric.
# Enter it with a leading apostrophe to force text, then re-apply the
# plain header format (bold/border/center, no quote-prefix) on top by
# copying it from an untouched header cell elsewhere in the workbook.
$ws.Range("B1").Value = "'2015"
$wsRef = $wb.Worksheets.Item("Potencia Acumulada - SIN (MW)")
$wsRef.Range("C1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A2").ClearFormats()
$ws.Range("A2").Value = "Expansão Centralizada"
$ws.Range("B2").Value = 573

$ws.Range("A3").ClearFormats()
$ws.Range("A3").Value = "Expansão por GD"
$ws.Range("B3").Value = 99
